$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EMPLOYEE DTR")

# Update REMARKS (column P) text for several rows, appending " R " and
# inserting "~ = " after "~OT " where applicable.
$ws.Range("P5").Value = '" ~OB Others|Roll-Out Pentstar Sports, Inc. | R "'
$ws.Range("P6").Value = '" ~OB Others|Roll-Out Pentstar Sports, Inc. | R "'
$ws.Range("P7").Value = '" ~OT ~ = Roll-Out Pentstar Sports, Inc. Robinsons Place Ermita ~OB Others|Roll-Out Pentstar Sports, Inc. | R "'
$ws.Range("P8").Value = '" ~OT ~ = Roll-Out Pilot Store SM MegaMall ~OB Others|Roll-Out Pentstar Sports, Inc. | R "'
$ws.Range("P9").Value = '" ~OT ~ = Roll-Out Pentstar Sports, Inc.  Newport, Paranque City ~OB Others|Roll-Out Pentstar Sports, Inc. | R "'
$ws.Range("P11").Value = '" ~OT ~ = Mall Integration, Pentstar Sports Inc. SM MegaMall"'
$ws.Range("P16").Value = '" ~OT ~ = GO-Live The Breakfast Club Pioneer Center, Pasig City"'

# Update overtime hours (column G) values
$ws.Range("G8").Value = 4.5
$ws.Range("G9").Value = 6.5
$ws.Range("G11").Value = 3.5
$ws.Range("G16").Value = 3.5
